$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.272.73"
$ws.Range("D3").Value = "1.619.55"
$ws.Range("E3").Value = "  +1.76%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'212.15"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'18.79"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.846.02"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "1.602.62"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +0.41%  "
$ws.Range("D15").Value = "'0.517"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "26.294.18"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "'62.41"
$ws.Range("E17").Value = "  +4.08%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'201.67"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").Value = "'9.33"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "'144.01"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").Value = "'15.15"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  +9.64%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("D32").Value = "'3.18"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.50"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.41"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").Value = "1.180.38"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("D38").Value = "'0.804"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'0.495"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").Value = "'5.34"
$ws.Range("E43").Value = "  +4.98%  "
$ws.Range("D44").Value = "1.756.62"
$ws.Range("E44").Value = "  +1.75%  "
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").Value = "0.0₆0104"
$ws.Range("E46").Value = "  +12.80%  "
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").Value = "'53.78"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "'0.410"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("E51").Value = "  -0.31%  "
